$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 219, shifting existing rows 219:260 down to 220:261
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new record
$ws.Cells.Item(219, 1).Value = 5
$ws.Cells.Item(219, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(219, 3).Value = "Maule"
$ws.Cells.Item(219, 4).Value = 44964
$ws.Cells.Item(219, 5).Value = 7
$ws.Cells.Item(219, 6).Value = 100112017
$ws.Cells.Item(219, 7).Value = "Apio"
$ws.Cells.Item(219, 8).Value = "Americana (o)"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 400
$ws.Cells.Item(219, 11).Value = 10000
$ws.Cells.Item(219, 12).Value = 10000
$ws.Cells.Item(219, 13).Value = 10000
$ws.Cells.Item(219, 14).Value = "$/docena de matas"
$ws.Cells.Item(219, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(219, 16).Value = 1667
$ws.Cells.Item(219, 17).Value = 6
$ws.Cells.Item(219, 18).Value = "Hortaliza"
